$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right after "2021-Q4" (before "总计"),
#    cloning the layout/formatting of the "2021-Q4" sheet and filling in the
#    new quarter's fund-holding figures.
# ---------------------------------------------------------------------------
$prevQuarter = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $prevQuarter)
$q1.Name = "2022-Q1"

# Clone cell formatting (styles/borders/fonts) from the previous quarter sheet
# (column A / row 1 is intentionally left out of the header copy since that
# sheet never populates A1).
$prevQuarter.Range("B1:H1").Copy()
$q1.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$prevQuarter.Range("A2:H3").Copy()
$q1.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "090019"
$q1.Range("C2").Value = "大成景恒混合A"
$q1.Range("D2").Value = "2.31"
$q1.Range("E2").Value = "93.51"
$q1.Range("F2").Value = "2.01"
$q1.Range("G2").Value = "0.0464"
$q1.Range("B2:G2").ClearFormats()
$q1.Range("H2").Value = 4

$q1.Range("A3").Value = 1
$q1.Range("B3:G3").NumberFormat = "@"
$q1.Range("B3").Value = "006038"
$q1.Range("C3").Value = "大成景恒混合C"
$q1.Range("D3").Value = "0.92"
$q1.Range("E3").Value = "93.51"
$q1.Range("F3").Value = "2.01"
$q1.Range("G3").Value = "0.0185"
$q1.Range("B3:G3").ClearFormats()
$q1.Range("H3").Value = 4

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: prepend a row for 2022-Q1, pushing the
#    existing 2021-Q4 / 2021-Q3 rows down.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.06

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# Restore the row-label style (s="2") on the new A2 cell to match its peers
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# Keep the originally active sheet selected (unchanged by this edit)
$wb.Worksheets.Item("2021-Q3").Activate()
